$wb = $excel.ActiveWorkbook

# "Forecast Comparison" sheet holds the per-week metrics, including the
# "Seasonality Index" column (L).
$wsForecast = $wb.Worksheets.Item("Forecast Comparison")

# New Seasonality Index values for rows 2-17 (column L)
$seasonalityIndex = @{
    2  = 0.95
    3  = 1.2
    4  = 0.97
    5  = 0.9
    6  = 0.89
    7  = 0.96
    8  = 1.09
    9  = 1.11
    10 = 1.07
    11 = 0.92
    12 = 1.2
    13 = 1.06
    14 = 1.09
    15 = 1.08
    16 = 1.09
    17 = 0.86
}

foreach ($row in $seasonalityIndex.Keys) {
    $wsForecast.Range("L$row").Value = $seasonalityIndex[$row]
}

# "Summary" sheet holds the aggregated totals, including the
# "Total Forecast (16 Weeks)" value in B9. This column stores its values
# as text, so prefix with an apostrophe to force Excel to keep "20" as a
# string rather than auto-converting it to a number.
$wsSummary = $wb.Worksheets.Item("Summary")
$wsSummary.Range("B9").Value = "'20"
